# "Generate Report for Handoff"
# Adds a new handoff-tracking row (for file 6911ee96-...md) to all three
# sheets: Overview, zh-cn, de-de. The new row follows the same layout /
# table structure as the existing row for the dc4713e4-...md file.

$wb = $excel.ActiveWorkbook

$mdNew   = "6911ee96-2ac0-4ec6-8f40-7805939eed31ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathNew = "e2e\6911ee96-2ac0-4ec6-8f40-7805939eed31ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$urlNew  = "https://github.com/OpenLocalizationTestOrg/oltest/blob/19f1ffe02e0e12a08deb59f88cb6af3fab2cc3e1/e2e/" + $mdNew

$zhCnXlfNew = "6911ee96-2ac0-4ec6-8f40-7805939eed31ooooooooooooooooooooooooooooooooooooooooooo.dcbc7e3d957f820167414ec3c9616bf595a6366c.zh-cn.xlf"
$deDeXlfNew = "6911ee96-2ac0-4ec6-8f40-7805939eed31ooooooooooooooooooooooooooooooooooooooooooo.dcbc7e3d957f820167414ec3c9616bf595a6366c.de-de.xlf"
$zhCnXlfOld = "dc4713e4-10c2-4544-b880-eb61f12efd1booooooooooooooooooooooooooooooooooooooooooo.85a8922fc61f66aa162a013dfd05725e41fcbc53.zh-cn.xlf"

$statusNew  = "Ready for handoff"
$dateGen    = "2016-08-12 02:50:24"
$dateHoXlf  = "2016-08-12 02:50:19"
$dateHandback = "2016-08-12 02:49:36"
$epoch      = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview": add row to its table, one entry per file (7 columns)
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$rowOv = $loOv.ListRows.Add()
$r = $loOv.DataBodyRange.Rows.Item($loOv.DataBodyRange.Rows.Count).Row

$wsOv.Range("A" + $r).Value = $mdNew
$wsOv.Range("B" + $r).Value = $pathNew
$wsOv.Hyperlinks.Add($wsOv.Range("B" + $r), $urlNew, "", "", $pathNew)
$wsOv.Range("C" + $r).Value = ".md"
$wsOv.Range("D" + $r).Value = ""
$wsOv.Range("E" + $r).Value = $statusNew
$wsOv.Range("F" + $r).Value = $statusNew
$wsOv.Range("G" + $r).Value = $dateGen
$wsOv.Range("G" + $r).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn": add matching row
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rz = $loZh.DataBodyRange.Rows.Item($loZh.DataBodyRange.Rows.Count).Row

$wsZh.Range("A" + $rz).Value = $mdNew
$wsZh.Hyperlinks.Add($wsZh.Range("A" + $rz), $urlNew, "", "", $mdNew)
$wsZh.Range("B" + $rz).Value = ".md"
$wsZh.Range("C" + $rz).Value = $statusNew
$wsZh.Range("D" + $rz).Value = "e2e"
$wsZh.Range("E" + $rz).Value = "ht"
$wsZh.Range("F" + $rz).Value = "'False"
$wsZh.Range("G" + $rz).Value = $epoch
$wsZh.Range("H" + $rz).Value = "'True"
$wsZh.Range("H" + $rz).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I" + $rz).Value = ""
$wsZh.Range("J" + $rz).Value = ""
$wsZh.Range("K" + $rz).Value = $zhCnXlfOld
$wsZh.Range("K" + $rz).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L" + $rz).Value = ""
$wsZh.Range("M" + $rz).Value = $dateHandback
$wsZh.Range("N" + $rz).Value = ""
$wsZh.Range("O" + $rz).Value = "'False"
$wsZh.Range("P" + $rz).Value = ""

# ---------------------------------------------------------------------
# Sheet "de-de": add matching row
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rd = $loDe.DataBodyRange.Rows.Item($loDe.DataBodyRange.Rows.Count).Row

$wsDe.Range("A" + $rd).Value = $mdNew
$wsDe.Hyperlinks.Add($wsDe.Range("A" + $rd), $urlNew, "", "", $mdNew)
$wsDe.Range("B" + $rd).Value = ".md"
$wsDe.Range("C" + $rd).Value = $statusNew
$wsDe.Range("D" + $rd).Value = "e2e"
$wsDe.Range("E" + $rd).Value = "ht"
$wsDe.Range("F" + $rd).Value = "'False"
$wsDe.Range("G" + $rd).Value = $deDeXlfNew
$wsDe.Range("H" + $rd).Value = $dateGen
$wsDe.Range("H" + $rd).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I" + $rd).Value = ""
$wsDe.Range("J" + $rd).Value = ""
$wsDe.Range("K" + $rd).Value = $zhCnXlfOld
$wsDe.Range("K" + $rd).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L" + $rd).Value = ""
$wsDe.Range("M" + $rd).Value = $dateHandback
$wsDe.Range("N" + $rd).Value = ""
$wsDe.Range("O" + $rd).Value = "'False"
$wsDe.Range("P" + $rd).Value = ""
